$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the term values (changing some terms of the matrix / constraint vector)
$ws.Range("E2").Value = -4
$ws.Range("F2").Value = -5
$ws.Range("K2").Value = -10

$ws.Range("E13").Value = 1
$ws.Range("F13").Value = 1

# Update the active cell selection to O8
$ws.Range("O8").Select()
